$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2" = 0.1508196721311476
    "C2" = 0.639344262295082
    "J2" = 0.01639344262295082
    "P2" = 0.1213114754098361
    "S2" = 0.07213114754098361
    "C3" = 0.0392156862745098
    "J3" = 0.009803921568627451
    "P3" = 0.7009803921568627
    "S3" = 0.25
    "J4" = 0.04545454545454546
    "P4" = 0.7045454545454546
    "S4" = 0.25
    "B6" = 0.08097165991902834
    "D6" = 0.01619433198380567
    "F6" = 0.08097165991902834
    "J6" = 0.2226720647773279
    "O6" = 0.008097165991902834
    "Q6" = 0.1740890688259109
    "R6" = 0.06072874493927125
    "S6" = 0.3562753036437247
    "B7" = 0.1442307692307692
    "D7" = 0.02884615384615385
    "F7" = 0.05288461538461538
    "J7" = 0.125
    "O7" = 0.009615384615384616
    "Q7" = 0.1923076923076923
    "R7" = 0.08653846153846154
    "S7" = 0.3605769230769231
    "B8" = 0.09239130434782608
    "D8" = 0.01449275362318841
    "F8" = 0.05072463768115942
    "J8" = 0.108695652173913
    "O8" = 0.02355072463768116
    "Q8" = 0.197463768115942
    "R8" = 0.06521739130434782
    "S8" = 0.447463768115942
    "B9" = 0.08994708994708994
    "D9" = 0.01587301587301587
    "F9" = 0.07407407407407407
    "J9" = 0.1058201058201058
    "O9" = 0.03174603174603174
    "Q9" = 0.164021164021164
    "R9" = 0.08994708994708994
    "S9" = 0.4285714285714285
    "B10" = 0.1133280127693535
    "D10" = 0.01915403032721468
    "E10" = 0.001596169193934557
    "F10" = 0.08379888268156424
    "J10" = 0.1037509976057462
    "O10" = 0.01436552274541101
    "Q10" = 0.209098164405427
    "R10" = 0.06943335993615324
    "S10" = 0.3854748603351955
    "G11" = 0.1219512195121951
    "J11" = 0.1310975609756098
    "K11" = 0.1829268292682927
    "L11" = 0.551829268292683
    "S11" = 0.01219512195121951
    "G12" = 0.7446808510638298
    "J12" = 0.1702127659574468
    "K12" = 0.005319148936170213
    "L12" = 0.03191489361702127
    "S12" = 0.04787234042553191
    "G13" = 0.6545454545454545
    "J13" = 0.3454545454545455
    "F15" = 0.004132231404958678
    "H15" = 0.1900826446280992
    "I15" = 0.07851239669421488
    "J15" = 0.3181818181818182
    "K15" = 0.1033057851239669
    "O15" = 0.05785123966942149
    "S15" = 0.2479338842975207
    "F16" = 0.02926829268292683
    "H16" = 0.175609756097561
    "I16" = 0.05853658536585366
    "J16" = 0.4
    "K16" = 0.1024390243902439
    "M16" = 0.01463414634146342
    "O16" = 0.06341463414634146
    "S16" = 0.1560975609756098
    "F17" = 0.01867219917012448
    "H17" = 0.2157676348547718
    "I17" = 0.07883817427385892
    "J17" = 0.3879668049792531
    "K17" = 0.08713692946058091
    "M17" = 0.01659751037344398
    "N17" = 0.002074688796680498
    "O17" = 0.05601659751037345
    "S17" = 0.1369294605809129
    "F18" = 0.04069767441860465
    "H18" = 0.2267441860465116
    "I18" = 0.1046511627906977
    "J18" = 0.3488372093023256
    "K18" = 0.09302325581395349
    "M18" = 0.02906976744186046
    "N18" = 0.005813953488372093
    "O18" = 0.04651162790697674
    "S18" = 0.1046511627906977
    "F19" = 0.01129146083274524
    "H19" = 0.2335920959774171
    "I19" = 0.07339449541284404
    "J19" = 0.3338038108680311
    "K19" = 0.1136203246294989
    "M19" = 0.02681721947776994
    "N19" = 0.0007057163020465773
    "O19" = 0.08045165843330981
    "S19" = 0.1263232180663373
}

foreach ($key in $updates.Keys) {
    $ws.Range($key).Value = $updates[$key]
}

Write-Output "Updated $($updates.Count) cells"